# Update national summary tables for the Transitions Rule with the newest
# airtoxics NATA data (allocation / transitions rule output refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "Means": Total Cancer Risk (row 9) and Total Respiratory (row 10) ---
$wsMeans = $wb.Worksheets.Item("Means")

$wsMeans.Range("B9").Value = 26
$wsMeans.Range("C9").Value = 23
$wsMeans.Range("D9").Value = 34
$wsMeans.Range("E9").Value = 31
$wsMeans.Range("F9").Value = 31
$wsMeans.Range("G9").Value = 31

$wsMeans.Range("B10").Value = 0.31
$wsMeans.Range("C10").Value = 0.27
$wsMeans.Range("D10").Value = 0.36
$wsMeans.Range("E10").Value = 0.36
$wsMeans.Range("F10").Value = 0.35
$wsMeans.Range("G10").Value = 0.36

# --- Sheet "Standard Deviations": Total Cancer Risk (row 9) and Total Respiratory (row 10) ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")

$wsSD.Range("B9").Value = 8.3
$wsSD.Range("C9").Value = 7.2
$wsSD.Range("D9").Value = 13
$wsSD.Range("E9").Value = 10
$wsSD.Range("F9").Value = 10
$wsSD.Range("G9").Value = 17

$wsSD.Range("B10").Value = 0.11
$wsSD.Range("C10").Value = 0.094
$wsSD.Range("D10").Value = 0.083
$wsSD.Range("E10").Value = 0.075
$wsSD.Range("F10").Value = 0.071
$wsSD.Range("G10").Value = 0.075
